$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = 25
$ws.Range("E15").Select()
